$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ServicesCategory")
$ws2 = $wb.Worksheets.Item("ServiceType")
$ws3 = $wb.Worksheets.Item("ServiceNameEBP")

# Sheet1 (ServicesCategory): row 3 is overwritten with what used to be row 4's
# data (column A explicitly becomes "testT4116" instead of the old row's
# "testT4116_Negative"); the old row 4 is then removed.
$ws1.Range("A3").Value = "testT4116"
$ws1.Range("B3").Value = 1
$ws1.Range("C3").Value = "NG1"
$ws1.Range("D3").Value = "click"
$ws1.Range("E3").Value = "n/a"
$ws1.Range("F3").Value = "click"
$ws1.Range("G3").Value = "autoText"
$ws1.Range("H3").Value = "n/a"
$ws1.Rows.Item(4).Delete()

# Sheet2 (ServiceType): same pattern, but column A keeps its original value.
$ws2.Range("A3").Value = "testT4116_Negative"
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = "NG1"
$ws2.Range("D3").Value = "click"
$ws2.Range("E3").Value = "click"
$ws2.Range("F3").Value = "n/a"
$ws2.Range("G3").Value = "click"
$ws2.Range("H3").Value = "autoText"
$ws2.Range("I3").Value = "n/a"
$ws2.Rows.Item(4).Delete()

# Sheet3 (ServiceNameEBP): same pattern; C3 keeps the quote-prefix ("text")
# cell style that the source row already carried, so it is entered with a
# leading apostrophe to reproduce that formatting.
$ws3.Range("A3").Value = "testT4116"
$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = "'NG1"
$ws3.Range("D3").Value = "click"
$ws3.Range("E3").Value = "n/a"
$ws3.Range("F3").Value = "click"
$ws3.Range("G3").Value = "autoText"
$ws3.Range("H3").Value = "No Evidence"
$ws3.Range("I3").Value = "Concerning Practice"
$ws3.Range("J3").Value = "n/a"
$ws3.Rows.Item(4).Delete()

# Selections / active sheet & cell to match the saved view state.
$ws2.Range("C9").Select() | Out-Null
$ws3.Range("E5").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("G4").Select() | Out-Null
